$d = $word.ActiveDocument

# Update AEA P&P report count: 572 -> 573 (whole word, so it doesn't touch "4572"-like numbers)
$d.Content.Find.Execute("572", $true, $true, $false, $false, $false,
                         $true, 1, $false, "573", 2)

# Update AEJ:Applied Economics report count: 973 -> 974 (whole word, so it doesn't touch "4973")
$d.Content.Find.Execute("973", $true, $true, $false, $false, $false,
                         $true, 1, $false, "974", 2)

# Update undergraduate internship paragraph
$d.Content.Find.Execute("have been engaged in this academic activity. Their names are listed in the Appendix. In addition, a pilot project in Summer 2024 provided interships to 9 undergraduates from various undergraduate institutions around the U.S.", $true, $false, $false, $false, $false,
                         $true, 1, $false, "have been engaged in this academic activity. Their names are listed in the Appendix. In addition, since the Summer of 2024, I have provided interships to 9 undergraduates each year from various undergraduate institutions around the U.S.", 2)

# Update graduate student count: six -> seven
$d.Content.Find.Execute("six graduate students", $true, $false, $false, $false, $false,
                         $true, 1, $false, "seven graduate students", 2)
